$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price/volume cells so values keep their exact
# literal representation (preserve trailing zeros, percent signs, etc.)
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6",
    "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10",
    "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14",
    "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18",
    "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23",
    "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38",
    "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42",
    "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46",
    "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50",
    "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row
# Row 2
$ws.Range("D2").Value = "307.19"
$ws.Range("E2").Value = "-0.14%"

# Row 3
$ws.Range("D3").Value = "40.50"
$ws.Range("E3").Value = "0.63%"

# Row 4
$ws.Range("D4").Value = "5.053"
$ws.Range("E4").Value = "0.00%"

# Row 5
$ws.Range("E5").Value = "-2.31%"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "1.601"
$ws.Range("E6").Value = "-2.51%"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9051"
$ws.Range("E7").Value = "-0.39%"

# Row 8
$ws.Range("D8").Value = "2.429"
$ws.Range("E8").Value = "-5.08%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1020"
$ws.Range("E9").Value = "-0.54%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1753"
$ws.Range("E10").Value = "0.18%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09091"
$ws.Range("E11").Value = "1.00%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04208"
$ws.Range("E12").Value = "-5.03%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "-0.51%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001256"
$ws.Range("E14").Value = "-1.72%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005818"
$ws.Range("E15").Value = "-2.67%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.353"
$ws.Range("E16").Value = "-0.33%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.264"
$ws.Range("E17").Value = "-1.38%"

# Row 18
$ws.Range("E18").Value = "-2.96%"

# Row 19
$ws.Range("D19").Value = "6.765"
$ws.Range("E19").Value = "-4.54%"

# Row 20
$ws.Range("D20").Value = "0.1366"
$ws.Range("E20").Value = "-1.43%"

# Row 21
$ws.Range("E21").Value = "2.48%"

# Row 22
$ws.Range("D22").Value = "0.04186"
$ws.Range("E22").Value = "0.37%"

# Row 23
$ws.Range("D23").Value = "0.001229"
$ws.Range("E23").Value = "1.23%"

# Row 24
$ws.Range("D24").Value = "0.004044"
$ws.Range("E24").Value = "-0.89%"

# Row 25
$ws.Range("D25").Value = "0.0001305"
$ws.Range("E25").Value = "6.14%"

# Row 26
$ws.Range("D26").Value = "0.0003013"
$ws.Range("E26").Value = "0.46%"

# Row 38
$ws.Range("D38").Value = "0.02386"
$ws.Range("E38").Value = "-0.72%"

# Row 39
$ws.Range("D39").Value = "0.05156"
$ws.Range("E39").Value = "-0.75%"

# Row 40
$ws.Range("D40").Value = "0.007758"
$ws.Range("E40").Value = "-2.92%"

# Row 41
$ws.Range("D41").Value = "0.1296"
$ws.Range("E41").Value = "-2.75%"

# Row 42
$ws.Range("D42").Value = "0.007062"
$ws.Range("E42").Value = "-5.81%"

# Row 43
$ws.Range("D43").Value = "0.001924"
$ws.Range("E43").Value = "-4.71%"

# Row 44
$ws.Range("D44").Value = "0.008509"
$ws.Range("E44").Value = "5.30%"

# Row 45
$ws.Range("D45").Value = "0.3350"
$ws.Range("E45").Value = "0.26%"

# Row 46
$ws.Range("D46").Value = "0.00006370"
$ws.Range("E46").Value = "-5.59%"

# Row 47
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "-0.47%"

# Row 48
$ws.Range("D48").Value = "0.004409"
$ws.Range("E48").Value = "6.80%"

# Row 49
$ws.Range("D49").Value = "0.006166"
$ws.Range("E49").Value = "85.68%"

# Row 50
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "-0.47%"

# Row 51
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "-0.47%"
